$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (pushes old row4 "Electrical Components" down to row5)
$ws.Rows.Item(4).Insert()

# Update row 2: "Electrical Equipment" -> "Electronic Components" (B2 stays "1")
$ws.Range("A2").Value = "Electronic Components"

# Update row 3: " Electronic Components" -> " Labor — Temporary", B3: 35 -> 1
$ws.Range("A3").Value = " Labor " + [char]0x2014 + " Temporary"
$ws.Range("B3").Value = "1"

# New row 4: "Electrical Equipment", 3
$ws.Range("A4").Value = "Electrical Equipment"
$ws.Range("B4").Value = "3"

# Row 5 (previously row 4): "Electrical Components" stays, B5: 37 -> 43
$ws.Range("B5").Value = "43"
